# Updated cryptos list — refreshed Coin/Link/Price/Volume(1h) snapshot.
#
# Column D ("Price") values are free-form text like "27.466.18" or
# "0.00001019" (locale-formatted strings, not real numbers — some even
# contain two '.' separators). Several of them otherwise parse as a
# valid Excel number (e.g. "1.005"), so a plain .Value assignment would
# silently coerce them into numeric cells and mangle values such as
# "5.270" -> 5.27. To keep them as genuine text (matching the original
# inline-string cells) we temporarily force a Text number format before
# writing, then restore the default "Normal" style so no stray
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: call this positionally (Set-TextValue $range $value) — passing the
# Range via a *named* parameter (-Range $range) loses the live COM binding
# in this runtime and silently no-ops, so positional args are required here.
function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  B = $null; C = $null; D = '27.466.18';   E = '  -0.42%  ' },
    @{ Row = 3;  B = $null; C = $null; D = '1.827.77';    E = '  -1.14%  ' },
    @{ Row = 4;  B = $null; C = $null; D = '1.005';       E = '  -0.77%  ' },
    @{ Row = 5;  B = $null; C = $null; D = '333.28';      E = '  -0.29%  ' },
    @{ Row = 6;  B = $null; C = $null; D = '1.004';       E = '  -0.65%  ' },
    @{ Row = 7;  B = $null; C = $null; D = '0.4586';      E = '  -0.30%  ' },
    @{ Row = 8;  B = $null; C = $null; D = '0.3823';      E = '  -1.76%  ' },
    @{ Row = 9;  B = $null; C = $null; D = '46.09';       E = '  +0.46%  ' },
    @{ Row = 10; B = $null; C = $null; D = '0.07828';     E = '  -1.17%  ' },
    @{ Row = 11; B = $null; C = $null; D = '0.9589';      E = '  -4.21%  ' },
    @{ Row = 12; B = $null; C = $null; D = '21.02';       E = '  -2.45%  ' },
    @{ Row = 13; B = $null; C = $null; D = '1.845.83';    E = '  -0.66%  ' },
    @{ Row = 14; B = $null; C = $null; D = '5.831';       E = '  -2.01%  ' },
    @{ Row = 15; B = $null; C = $null; D = '7.053';       E = '  -1.62%  ' },
    @{ Row = 16; B = $null; C = $null; D = '1.006';       E = '  -0.68%  ' },
    @{ Row = 17; B = $null; C = $null; D = $null;         E = '  +1.21%  ' },
    @{ Row = 18; B = $null; C = $null; D = '0.06575';     E = '  -1.81%  ' },
    @{ Row = 19; B = $null; C = $null; D = '0.00001019';  E = '  -1.57%  ' },
    @{ Row = 20; B = $null; C = $null; D = '17.07';       E = '  -0.61%  ' },
    @{ Row = 21; B = $null; C = $null; D = $null;         E = '  -0.56%  ' },
    @{ Row = 22; B = $null; C = $null; D = '27.473.98';   E = '  -0.45%  ' },
    @{ Row = 23; B = $null; C = $null; D = '5.288';       E = '  -2.29%  ' },
    @{ Row = 24; B = $null; C = $null; D = '10.77';       E = '  -1.29%  ' },
    @{ Row = 25; B = $null; C = $null; D = '2.274';       E = '  -1.35%  ' },
    @{ Row = 26; B = $null; C = $null; D = '2.067.83';    E = '  -0.40%  ' },
    @{ Row = 27; B = $null; C = $null; D = '158.97';      E = '  -0.07%  ' },
    @{ Row = 28; B = $null; C = $null; D = '19.33';       E = '  -1.10%  ' },
    @{ Row = 29; B = $null; C = $null; D = '2.033';       E = '  -4.51%  ' },
    @{ Row = 30; B = $null; C = $null; D = '5.270';       E = '  -3.14%  ' },
    @{ Row = 31; B = $null; C = $null; D = '117.62';      E = '  -3.09%  ' },
    @{ Row = 32; B = $null; C = $null; D = '0.09386';     E = '  -0.18%  ' },
    @{ Row = 33; B = $null; C = $null; D = '0.9254';      E = '  -4.81%  ' },
    @{ Row = 34; B = $null; C = $null; D = '3.575';       E = '  -1.23%  ' },
    @{ Row = 35; B = $null; C = $null; D = '5.195';       E = '  -2.07%  ' },
    @{ Row = 36; B = $null; C = $null; D = '1.313';       E = '  -1.54%  ' },
    @{ Row = 37; B = $null; C = $null; D = '0.05919';     E = '  -1.37%  ' },
    @{ Row = 38; B = $null; C = $null; D = '0.02183';     E = '  -2.08%  ' },
    @{ Row = 39; B = $null; C = $null; D = '8.095';       E = '  -2.95%  ' },
    @{ Row = 40; B = $null; C = $null; D = $null;         E = '  -0.58%  ' },
    @{ Row = 41; B = $null; C = $null; D = '1.142';       E = '  -3.66%  ' },
    @{ Row = 42; B = $null; C = $null; D = '0.5719';      E = '  -3.30%  ' },
    @{ Row = 43; B = $null; C = $null; D = $null;         E = '  -2.29%  ' },
    @{ Row = 44; B = $null; C = $null; D = '9.925';       E = '  -4.72%  ' },
    @{ Row = 45; B = $null; C = $null; D = '1.268';       E = '  +1.95%  ' },
    @{ Row = 46; B = 'EnergySwap';   C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens';      D = '11.78';  E = '  -3.13%  ' },
    @{ Row = 47; B = 'Decentraland'; C = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D = '0.5369'; E = '  -3.77%  ' },
    @{ Row = 48; B = $null; C = $null; D = '1.897';       E = '  -0.57%  ' },
    @{ Row = 49; B = $null; C = $null; D = '0.06821';     E = '  +1.86%  ' },
    @{ Row = 50; B = $null; C = $null; D = '110.08';      E = '  -0.96%  ' },
    @{ Row = 51; B = $null; C = $null; D = $null;         E = '  -32.65%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.B) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($null -ne $u.C) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($null -ne $u.D) {
        Set-TextValue $ws.Range("D$row") $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
